$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.745220526266252
$ws.Range("C2").Value = 0.0919792231427196
$ws.Range("D2").Value = -8.10205284197639

$ws.Range("B3").Value = -0.467632222796714
$ws.Range("C3").Value = 0.0919925416445528
$ws.Range("D3").Value = -5.08337104766149

$ws.Range("B4").Value = -0.00839567998234618
$ws.Range("C4").Value = 0.0076652298887339
$ws.Range("D4").Value = -1.09529395780887

$ws.Range("B5").Value = 0.409586883899209

$ws.Range("B6").Value = 0.165685936536102
